$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Save" header in H1, matching the formatting (style) of the other
# header cells (e.g. G1 = "sum") by copying formats across.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Populate the new "Save" column data (H2:H9).
$saveValues = @(0, 0, 0, 1, 0, 0, 1, 0)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
